# Auto-generated Excel COM-interop script applying scheduled-runner updates
# to the Sheets workbook's per-leve market price / profit columns (H:N).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4165.241
$ws.Range("J64").Value = 4192.1924
$ws.Range("L64").Value = 4192.1924
$ws.Range("N64").Value = -4688.1924

$ws.Range("H67").Value = 4165.241
$ws.Range("J67").Value = 4192.1924
$ws.Range("L67").Value = 4192.1924
$ws.Range("N67").Value = -5908.1924

$ws.Range("H70").Value = 2287.76
$ws.Range("I70").Value = 2109.6667
$ws.Range("J70").Value = 2387.9375
$ws.Range("K70").Value = 6329.000100000001
$ws.Range("L70").Value = 7163.8125
$ws.Range("M70").Value = -6059.000100000001
$ws.Range("N70").Value = -7703.8125

$ws.Range("H73").Value = 2287.76
$ws.Range("I73").Value = 2109.6667
$ws.Range("J73").Value = 2387.9375
$ws.Range("K73").Value = 6329.000100000001
$ws.Range("L73").Value = 7163.8125
$ws.Range("M73").Value = -5393.000100000001
$ws.Range("N73").Value = -9035.8125

$ws.Range("H87").Value = 92672
$ws.Range("J87").Value = 92672
$ws.Range("L87").Value = 92672
$ws.Range("N87").Value = -95168

$ws.Range("H90").Value = 92672
$ws.Range("J90").Value = 92672
$ws.Range("L90").Value = 278016
$ws.Range("N90").Value = -290496

$ws.Range("H133").Value = 86499.75
$ws.Range("J133").Value = 86499.75
$ws.Range("L133").Value = 86499.75
$ws.Range("N133").Value = -96619.75

$ws.Range("H138").Value = 2522.8777
$ws.Range("J138").Value = 3280.3508
$ws.Range("L138").Value = 9841.0524
$ws.Range("N138").Value = -20121.0524

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 2474.5
$ws.Range("J4").Value = 2998
$ws.Range("L4").Value = 2998
$ws.Range("N4").Value = -3230

$ws.Range("H32").Value = 9807338
$ws.Range("I32").Value = 10640858
$ws.Range("J32").Value = 13474.5
$ws.Range("K32").Value = 10640858
$ws.Range("L32").Value = 13474.5
$ws.Range("M32").Value = -10640571
$ws.Range("N32").Value = -14048.5

$ws.Range("H110").Value = 1124.5238
$ws.Range("I110").Value = 1149.3158
$ws.Range("J110").Value = 889
$ws.Range("K110").Value = 1149.3158
$ws.Range("L110").Value = 889
$ws.Range("M110").Value = 895.6841999999999
$ws.Range("N110").Value = -4979

$ws.Range("H122").Value = 1805.5
$ws.Range("I122").Value = 1805.5
$ws.Range("K122").Value = 5416.5
$ws.Range("M122").Value = -2966.5

$ws.Range("H132").Value = 6173.1113
$ws.Range("I132").Value = 2139.9
$ws.Range("J132").Value = 11214.625
$ws.Range("K132").Value = 6419.700000000001
$ws.Range("L132").Value = 33643.875
$ws.Range("M132").Value = -3889.700000000001
$ws.Range("N132").Value = -38703.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H33").Value = 8000
$ws.Range("I33").Value = 8000
$ws.Range("K33").Value = 8000
$ws.Range("M33").Value = -7664

$ws.Range("H105").Value = 2657.4375
$ws.Range("I105").Value = 2539.6667
$ws.Range("K105").Value = 2539.6667
$ws.Range("M105").Value = -792.6667000000002

$ws.Range("H107").Value = 1648.5333
$ws.Range("I107").Value = 1181.2222
$ws.Range("J107").Value = 5854.3335
$ws.Range("K107").Value = 1181.2222
$ws.Range("L107").Value = 5854.3335
$ws.Range("M107").Value = 738.7778000000001
$ws.Range("N107").Value = -9694.333500000001

$ws.Range("H134").Value = 24037.092
$ws.Range("I134").Value = 1105.35
$ws.Range("K134").Value = 3316.05
$ws.Range("M134").Value = -781.0499999999997

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 788585.4399999999
$ws.Range("I31").Value = 17392.363
$ws.Range("K31").Value = 17392.363
$ws.Range("M31").Value = -17097.363

$ws.Range("H34").Value = 788585.4399999999
$ws.Range("I34").Value = 17392.363
$ws.Range("K34").Value = 17392.363
$ws.Range("M34").Value = -17190.363

$ws.Range("H107").Value = 861.3333
$ws.Range("I107").Value = 804.4
$ws.Range("K107").Value = 804.4
$ws.Range("M107").Value = 1115.6

$ws.Range("H125").Value = 45285.715
$ws.Range("J125").Value = 45285.715
$ws.Range("L125").Value = 45285.715
$ws.Range("N125").Value = -50205.715

$ws.Range("H134").Value = 233960.53
$ws.Range("I134").Value = 251375.9
$ws.Range("K134").Value = 754127.7
$ws.Range("M134").Value = -751592.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 1542.7778
$ws.Range("I12").Value = 3113.7144
$ws.Range("J12").Value = 543.0909
$ws.Range("K12").Value = 9341.143199999999
$ws.Range("L12").Value = 1629.2727
$ws.Range("M12").Value = -9168.143199999999
$ws.Range("N12").Value = -1975.2727

$ws.Range("H50").Value = 454.04166
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 454.04166
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 1362.12498
$ws.Range("M50").Value = $null
$ws.Range("N50").Value = -2324.12498

$ws.Range("H51").Value = 24835.5
$ws.Range("I51").Value = 16666
$ws.Range("K51").Value = 49998
$ws.Range("M51").Value = -49538

$ws.Range("H53").Value = 454.04166
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 454.04166
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 1362.12498
$ws.Range("M53").Value = $null
$ws.Range("N53").Value = -2324.12498

$ws.Range("H68").Value = 202531.8
$ws.Range("I68").Value = 1874.5
$ws.Range("K68").Value = 5623.5
$ws.Range("M68").Value = -4812.5

$ws.Range("H71").Value = 202531.8
$ws.Range("I71").Value = 1874.5
$ws.Range("K71").Value = 16870.5
$ws.Range("M71").Value = -12814.5

$ws.Range("H134").Value = 7935.6
$ws.Range("I134").Value = 2888
$ws.Range("J134").Value = 9197.5
$ws.Range("K134").Value = 8664
$ws.Range("L134").Value = 27592.5
$ws.Range("M134").Value = -3594
$ws.Range("N134").Value = -37732.5

$ws.Range("H139").Value = 3047.45
$ws.Range("J139").Value = 2956.4285
$ws.Range("L139").Value = 8869.2855
$ws.Range("N139").Value = -19149.2855

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 590
$ws.Range("I2").Value = 307.7143
$ws.Range("J2").Value = 872.2857
$ws.Range("K2").Value = 307.7143
$ws.Range("L2").Value = 872.2857
$ws.Range("M2").Value = -194.7143
$ws.Range("N2").Value = -1098.2857

$ws.Range("H122").Value = 2461
$ws.Range("I122").Value = 2149.5
$ws.Range("J122").Value = 2668.6667
$ws.Range("K122").Value = 6448.5
$ws.Range("L122").Value = 8006.000100000001
$ws.Range("M122").Value = -3998.5
$ws.Range("N122").Value = -12906.0001

$ws.Range("H132").Value = 76940150
$ws.Range("I132").Value = 90911090
$ws.Range("J132").Value = 100013.5
$ws.Range("K132").Value = 272733270
$ws.Range("L132").Value = 300040.5
$ws.Range("M132").Value = -272730740
$ws.Range("N132").Value = -305100.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 62668.65
$ws.Range("I7").Value = 2499.875
$ws.Range("K7").Value = 2499.875
$ws.Range("M7").Value = -2387.875

$ws.Range("H45").Value = 13750
$ws.Range("I45").Value = 13750
$ws.Range("K45").Value = 13750
$ws.Range("M45").Value = -13343

$ws.Range("H61").Value = 1526.6364
$ws.Range("I61").Value = 1489.3
$ws.Range("K61").Value = 1489.3
$ws.Range("M61").Value = -1287.3

$ws.Range("H113").Value = 1526.6364
$ws.Range("I113").Value = 1489.3
$ws.Range("K113").Value = 1489.3
$ws.Range("M113").Value = 680.7

$ws.Range("H126").Value = 62668.65
$ws.Range("I126").Value = 2499.875
$ws.Range("K126").Value = 7499.625
$ws.Range("M126").Value = -5029.625

$ws.Range("H132").Value = 34384.43
$ws.Range("I132").Value = 7027.5415
$ws.Range("K132").Value = 21082.6245
$ws.Range("M132").Value = -18552.6245

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 962.5
$ws.Range("I113").Value = 950
$ws.Range("K113").Value = 2850
$ws.Range("M113").Value = -680

$ws.Range("H125").Value = 96497.25
$ws.Range("J125").Value = 96497.25
$ws.Range("L125").Value = 96497.25
$ws.Range("N125").Value = -106337.25

$ws.Range("H132").Value = 2854.8333
$ws.Range("I132").Value = 2336.75
$ws.Range("K132").Value = 7010.25
$ws.Range("M132").Value = -4480.25
